# Updates cryptos list (Coin, Link, Price, Volume(1h)) to latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='68.487.91'; E='  +0.91%  '},
    @{Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='3.815.71'; E='  -0.14%  '},
    @{Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='0.992'; E='  -0.77%  '},
    @{Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='611.84'; E='  +1.08%  '},
    @{Row=6; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='163.91'; E='  -1.30%  '},
    @{Row=7; B='LidoStakedEther'; C='https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'; D='3.814.85'; E='  -0.06%  '},
    @{Row=8; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='1.00'; E='  -0.02%  '},
    @{Row=9; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.517'; E='  -0.38%  '},
    @{Row=10; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.160'; E='  -0.03%  '},
    @{Row=11; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.450'; E='  -0.64%  '},
    @{Row=12; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='6.79'; E='  +7.29%  '},
    @{Row=13; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.0000247'; E='  -2.08%  '},
    @{Row=14; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='35.19'; E='  -2.62%  '},
    @{Row=15; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='4.455.53'; E='  -0.05%  '},
    @{Row=16; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='3.892.62'; E='  +1.96%  '},
    @{Row=17; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='68.525.04'; E='  +0.95%  '},
    @{Row=18; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='18.14'; E='  -1.41%  '},
    @{Row=19; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='7.10'; E='  -0.11%  '},
    @{Row=20; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.113'; E='  -0.24%  '},
    @{Row=21; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='462.76'; E='  -0.31%  '},
    @{Row=22; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='9.63'; E='  -2.59%  '},
    @{Row=23; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.700'; E='  -0.50%  '},
    @{Row=24; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='83.47'; E='  +0.05%  '},
    @{Row=25; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.0000146'; E='  -1.09%  '},
    @{Row=26; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='12.00'; E='  -1.20%  '},
    @{Row=27; B='Fetch.AI'; C='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D='2.11'; E='  -0.90%  '},
    @{Row=28; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.00'; E='  +0.08%  '},
    @{Row=29; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='9.99'; E='  -0.49%  '},
    @{Row=30; B='WrappedeETH'; C='https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'; D='3.957.69'; E='  -0.32%  '},
    @{Row=31; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='2.63'; E='  -6.05%  '},
    @{Row=32; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='2.20'; E='  -1.04%  '},
    @{Row=33; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='7.24'; E='  -2.63%  '},
    @{Row=34; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='28.93'; E='  -2.39%  '},
    @{Row=35; B='Binance-PegBSC-USD'; C='https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; D='1.00'; E='  -0.11%  '},
    @{Row=36; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='9.08'; E='  -0.51%  '},
    @{Row=37; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.101'; E='  +0.29%  '},
    @{Row=38; B='Kaspa'; C='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D='0.146'; E='  +5.49%  '},
    @{Row=39; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='5.89'; E='  +0.92%  '},
    @{Row=40; B='Mantle'; C='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D='0.980'; E='  -1.79%  '},
    @{Row=41; B='FirstDigitalUSD'; C='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D='1.00'; E='  +0.02%  '},
    @{Row=42; B='dogwifhat'; C='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D='3.12'; E='  -3.74%  '},
    @{Row=43; B='USDe'; C='https://coinranking.com/coin/exbfr2U-0+usde-usde'; D='1.00'; E='  +0.04%  '},
    @{Row=44; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='153.88'; E='  +1.49%  '},
    @{Row=45; B='Arweave'; C='https://coinranking.com/coin/7XWg41D1+arweave-ar'; D='43.10'; E='  -4.25%  '},
    @{Row=46; B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='0.296'; E='  -1.62%  '},
    @{Row=47; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='46.60'; E='  -2.47%  '},
    @{Row=48; B='ONDO'; C='https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'; D='1.39'; E='  +0.12%  '},
    @{Row=49; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='8.37'; E='  -0.11%  '},
    @{Row=50; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='1.86'; E='  +0.15%  '},
    @{Row=51; B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='379.32'; E='  -2.99%  '}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    # Price column holds numeric-looking text (e.g. "68.487.91", "0.992");
    # force text format so Excel doesn't reinterpret it as a number/date.
    $dCell = $ws.Cells.Item($r.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $r.D
    # Volume(1h) already contains spaces/percent sign, so it is kept as text naturally.
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
